# Apply crypto price/volume updates to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '70.385.59'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.622.08'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +2.75%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '602.01'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.74%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '196.85'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.27%  '
$ws.Cells.Item(7, 5).Value = '  -0.90%  '
$ws.Cells.Item(9, 5).Value = '  +6.54%  '
$ws.Cells.Item(10, 5).Value = '  -0.53%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '53.25'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -1.04%  '
$ws.Cells.Item(12, 5).Value = '  +0.71%  '
$ws.Cells.Item(13, 5).Value = '  +0.32%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.194.87'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.71%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '605.51'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.06%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.01'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.92%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '70.448.18'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.61%  '
$ws.Cells.Item(18, 2).Value = 'Chainlink'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.03'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.40%  '
$ws.Cells.Item(19, 2).Value = 'WrappedEther'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.600.87'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.90%  '
$ws.Cells.Item(20, 5).Value = '  +1.35%  '
$ws.Cells.Item(21, 5).Value = '  +0.47%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.08'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -1.36%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.23'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.32%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '103.26'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.22%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.60'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.89%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.98'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -6.63%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.62'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -2.54%  '
$ws.Cells.Item(28, 5).Value = '  +0.76%  '
$ws.Cells.Item(29, 5).Value = '  +1.20%  '
$ws.Cells.Item(30, 5).Value = '  +8.32%  '
$ws.Cells.Item(31, 5).Value = '  +2.47%  '
$ws.Cells.Item(32, 5).Value = '  -1.50%  '
$ws.Cells.Item(33, 5).Value = '  +0.76%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '63.30'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.29%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0888'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.12%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.937.17'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +5.59%  '
$ws.Cells.Item(37, 2).Value = 'Dai'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.03%  '
$ws.Cells.Item(38, 2).Value = 'Bittensor'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '521.45'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +6.79%  '
$ws.Cells.Item(39, 5).Value = '  -0.52%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '36.68'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.01%  '
$ws.Cells.Item(41, 5).Value = '  -1.35%  '
$ws.Cells.Item(42, 5).Value = '  -2.48%  '
$ws.Cells.Item(43, 5).Value = '  +2.01%  '
$ws.Cells.Item(44, 5).Value = '  +1.45%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.52'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +6.55%  '
$ws.Cells.Item(46, 5).Value = '  +1.80%  '
$ws.Cells.Item(47, 5).Value = '  -0.24%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.56'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.37%  '
$ws.Cells.Item(49, 5).Value = '  -0.25%  '
$ws.Cells.Item(50, 5).Value = '  -0.22%  '
$ws.Cells.Item(51, 5).Value = '  +0.48%  '

Write-Host "Applied" 81 "cell updates."
